$wb = $excel.ActiveWorkbook

# --- References to existing sheets ---
$totalSheet = $wb.Worksheets.Item("总计")        # currently sheetId=5, rId5, name "总计"
$q4Sheet    = $wb.Worksheets.Item("2021-Q4")     # header layout to reuse for the new "2022-Q1" sheet

# --- 1) Duplicate the current "总计" sheet (placed right after it) so the ---
#        copy inherits all of its formatting (header style, column A style, ---
#        page margins, etc.) for free. The duplicate becomes the NEW "总计" ---
#        sheet (sheetId=6), while the original sheet keeps its identity ---
#        (sheetId=5) and is repurposed into "2022-Q1". ---
$totalSheet.Copy($null, $totalSheet)
$newTotal = $wb.Worksheets.Item($totalSheet.Index + 1)

# --- 2) Rename the sheets into their final places. ---
$totalSheet.Name = "2022-Q1"
$newTotal.Name = "总计"

# --- 3) Rebuild the new "总计" totals table (adds the 2022-Q1 row on top, ---
#        pushing the rest down by one row). The header (row 1) and A2:A5 ---
#        style already match (inherited from the copy); row 6 is new, so ---
#        stamp it with the same column-A style before clearing + rewriting ---
#        the data cells. ---
$newTotal.Range("A5").Copy()
$newTotal.Range("A6").PasteSpecial(-4122)

$newTotal.Range("A2:D6").ClearContents()

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0.08

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 5
$newTotal.Range("D3").Value = 0.7

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q2"
$newTotal.Range("C4").Value = 9
$newTotal.Range("D4").Value = 1.77

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q1"
$newTotal.Range("C5").Value = 6
$newTotal.Range("D5").Value = 1.31

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2020-Q4"
$newTotal.Range("C6").Value = 8
$newTotal.Range("D6").Value = 1.06

# --- 4) Turn the repurposed sheet into the new "2022-Q1" fund-holdings ---
#        table: wipe the old totals content and reuse the fund-table ---
#        header/column formatting from the "2021-Q4" sheet. ---
$totalSheet.Cells.Clear()

$q4Sheet.Range("B1:H1").Copy()
$totalSheet.Range("B1:H1").PasteSpecial(-4122)

$q4Sheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

# B2:G2 hold text look-alikes of numbers (fund code / ratios) in the source
# data, so force a text number-format before writing them to avoid Excel
# auto-coercing them into real numbers (which would also strip the leading
# zero from the fund code). Restore the (unstyled) "Normal" cell style right
# afterwards so the stored cell doesn't pick up a stray style index.
$totalSheet.Range("B2:G2").NumberFormat = "@"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "005433"
$totalSheet.Range("C2").Value = "申万菱信医药先锋股票"
$totalSheet.Range("D2").Value = "2.20"
$totalSheet.Range("E2").Value = "90.81"
$totalSheet.Range("F2").Value = "3.69"
$totalSheet.Range("G2").Value = "0.0812"
$totalSheet.Range("H2").Value = 8

$totalSheet.Range("B2:G2").Style = "Normal"
